# Apply the edit described in the diff:
# - On "Hoja2" worksheet, add a new row 13 with a new value in B13:
#   "Actividad de prueba excel" (this becomes a new shared string).
# - The worksheet's used range / dimension grows from B4:B12 to B4:B13.
# - The active sheet selection on Hoja2 moves to M13.

$wb = $excel.ActiveWorkbook

# Locate "Hoja2" (the sheet that receives the new activity row).
$ws = $wb.Worksheets.Item("Hoja2")

# Write the new activity text into B13.
$ws.Range("B13").Value = "Actividad de prueba excel"

# Make this worksheet active and move the selection to M13, matching the
# final selection recorded in the workbook (activeCell="M13" sqref="M13").
$ws.Activate()
$ws.Range("M13").Select()
